$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D = "last_edited_time". Most rows that previously held
# 2024-08-30T20:17:00.000Z get bumped to 2024-08-31T05:43:00.000Z,
# except row 7 which gets its own distinct timestamp 2024-08-31T05:40:00.000Z.
$rowsCommon = @(2,3,5,6,8,11,13)
foreach ($r in $rowsCommon) {
    $ws.Cells.Item($r, 4).Value = "2024-08-31T05:43:00.000Z"
}

$ws.Cells.Item(7, 4).Value = "2024-08-31T05:40:00.000Z"
